$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data per the Tue May  7 14:22:42 UTC 2024
# GitHub Actions refresh. Rows 30 and 31 swap their Coin/Link/Price/Volume
# content (ImmutableX now ranks above PancakeSwap).
# NumberFormat is forced to Text ("@") immediately before writing any Price
# value that looks like a plain decimal number, so Excel keeps storing it as
# a text string (matching the workbook's inlineStr cells) instead of
# auto-converting it to a floating point number.

$ws.Range('D2').Value = '63.399.37'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '3.066.47'
$ws.Range('E3').Value = '  -1.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.58'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.45'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '3.066.32'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.93'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('E12').Value = '  -2.38%  '
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.55'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '3.570.44'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.20'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '63.351.61'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '3.063.87'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '486.75'
$ws.Range('E20').Value = '  +2.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.52'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.707'
$ws.Range('E22').Value = '  -4.48%  '
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.40'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.23'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.80'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('E27').Value = '  +5.45%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.49'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.33'
$ws.Range('E33').Value = '  -1.04%  '
$ws.Range('E34').Value = '  -4.58%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').Value = '0.0₃0824'
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('E38').Value = '  -5.56%  '
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.74'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '440.18'
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  +2.17%  '
$ws.Range('E45').Value = '  -2.72%  '
$ws.Range('D46').Value = '2.821.49'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.59'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.00'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.41'
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E51').Value = '  -2.65%  '

Write-Host "Updated cryptos list"
